# Add a new "AIOwnerID" column to the NPC sheet, inserted immediately
# before the existing "AI"/DescID column (old column AI, now shifted to AJ).
#
# This mirrors an Excel "Insert Column" at AI: the previously-existing
# columns AI and AJ shift right to AJ and AK respectively (their values,
# styles and the data validation / dimension / column-width metadata move
# with them automatically). We then populate the freshly inserted AI
# column:
#   - row 1  (header)      -> "AIOwnerID"
#   - rows 2-9 (meta rows: type/object flags) -> copied from column AH
#     (the sibling "MasterID" owner-id column), since the new field shares
#     the same type/flag metadata
#   - row 10 (description)  -> "AI"
#   - rows 11-62 (data)     -> 0

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert the new column; everything at/after AI shifts right by one.
$ws.Columns("AI").Insert()

# Row 1: new column header.
$ws.Range("AI1").Value = "AIOwnerID"

# Rows 2-9: copy the meta/flag values from column AH (the MasterID column)
for ($r = 2; $r -le 9; $r++) {
    $srcCell = $ws.Cells.Item($r, 34)   # column AH
    $dstCell = $ws.Cells.Item($r, 35)   # column AI
    $dstCell.Value = $srcCell.Value2
}

# Row 10: short description for the new field.
$ws.Range("AI10").Value = "AI"

# Rows 11-62: data values, all zero.
for ($r = 11; $r -le 62; $r++) {
    $ws.Cells.Item($r, 35).Value = 0
}

# Restore the column width on the new AI column to match its neighbour AH
# (Excel normally carries formatting over from the column to the left on
# insert).
$ws.Columns("AI").ColumnWidth = $ws.Columns("AH").ColumnWidth

# Update the selection to reflect where editing left off.
$ws.Range("AI11").Select() | Out-Null
